$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("E9").Value = "  +4.98%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E14").Value = "  -5.35%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  -4.10%  "
$ws.Range("E51").Value = "  +0.00%  "
